# "Đã cập nhập trạng thái" — widen the "Trạng thái" / "Ghi chú" columns of the
# plan table and mark the finished tasks as "Đã hoàn thành".

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Resize the table and its last two columns (values are OOXML twips;
#     Word's object model works in points, so divide by 20). ---
$t.PreferredWidth = 14317 / 20.0        # w:tblW  14038 -> 14317
$t.Columns.Item(7).Width = 1626 / 20.0  # "Trạng thái" column  1317 -> 1626
$t.Columns.Item(8).Width = 992 / 20.0   # last column          1022 -> 992

# --- Mark the rows whose work is complete. ---
# Row 1 is the header; data rows whose "Trạng thái" cell gets the new
# "Đã hoàn thành" note are the ones below (table has 36 rows total).
$doneRows = @(3, 4, 5, 6, 7, 8, 9, 10, 12, 13, 14, 15, 17, 18)

foreach ($rowNum in $doneRows) {
    $cell = $t.Cell($rowNum, 7)
    $cellRange = $cell.Range
    # Position right before the cell's trailing paragraph mark so the new
    # text lands at the end of the existing content instead of replacing it.
    $insertAt = $cellRange.End - 1
    $insertionPoint = $d.Range($insertAt, $insertAt)
    $insertionPoint.InsertAfter("Đã hoàn thành")
}
